# Daily scrape update - 2025-12-17 03:30:24 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data set (opportunity rows 2-13) ---------------------------------
# Columns: A=ID  B=LINK  C=TITLE  D=COUNTRY  E=PREMIUM  F=APPLICANTS  G=DURATION  H=ORGANIZATION
$rows = @(
    @("1330655", "https://aiesec.org/opportunity/global-talent/1330655", "Software Engineer Ruby on Rails/React/Flutter/AI/ML", "Lahore, Pakistan", "No", "0 applicants", "9 - 12 Weeks", "Arkhitech"),
    @("1330654", "https://aiesec.org/opportunity/global-talent/1330654", "UI/UX Designer", "Lahore, Pakistan", "No", "0 applicants", "9 - 12 Weeks", "Arkhitech"),
    @("1330653", "https://aiesec.org/opportunity/global-talent/1330653", "Business Development Manager", "Lahore, Pakistan", "No", "0 applicants", "9 - 12 Weeks", "Arkhitech"),
    @("1330646", "https://aiesec.org/opportunity/global-talent/1330646", "Software Developer", "Μαρούσι, Ελλάδα", "No", "1 applicant", "9 - 12 Weeks", "Inventio"),
    @("1330644", "https://aiesec.org/opportunity/global-talent/1330644", "Business Analyst", "Μαρούσι, Ελλάδα", "No", "2 applicants", "9 - 12 Weeks", "Inventio"),
    @("1330581", "https://aiesec.org/opportunity/global-talent/1330581", "Procurement Intern", "Panamá, Provincia de Panamá, Panamá", "No", "5 applicants", "6 - 18 Months", "Red Bull Panamá"),
    @("1328566", "https://aiesec.org/opportunity/global-talent/1328566", "HR Intern", "Santiago, Región Metropolitana, Chile", "No", "108 applicants", "6 - 18 Months", "Boehringer Ingelheim in Chile"),
    @("1328442", "https://aiesec.org/opportunity/global-talent/1328442", "Brand Ambassador", "台灣臺北", "No", "97 applicants", "3 - 6 Months", "Din Tai Fung Restaurant Co., Ltd."),
    @("1328023", "https://aiesec.org/opportunity/global-talent/1328023", "Marketing - Intern", "Nugegoda, Sri Lanka", "No", "36 applicants", "3 - 6 Months", "Raffles Consolidated Pvt Ltd"),
    @("1327811", "https://aiesec.org/opportunity/global-talent/1327811", "Software Engineering Intern", "Colombo, Sri Lanka", "No", "122 applicants", "3 - 6 Months", "Envision Circle (Pvt) Ltd"),
    @("1327778", "https://aiesec.org/opportunity/global-talent/1327778", "Digital Content & Stakeholder Engagement Intern", "Colombo, Sri Lanka", "No", "19 applicants", "6 - 18 Months", "Solutions Ground (Pvt) Ltd"),
    @("1327500", "https://aiesec.org/opportunity/global-talent/1327500", "Intern/Junior Engineer", "Barendrecht, Nederland", "No", "153 applicants", "6 - 18 Months", "CTS Offshore and Marine BV")
)

# Opportunity IDs in column A are purely numeric strings; force the cells to
# Text so they keep serializing like the original scrape (no numeric coercion).
$ws.Range("A2:A13").NumberFormat = "@"

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}

# Row 2's PREMIUM cell no longer needs the highlighted (yellow) style, since
# it is no longer "Yes" -- restore the default style like the rest of the column.
$ws.Range("E2").Style = "Normal"

# --- Column width tweaks ---------------------------------------------------
# (Excel's ColumnWidth is offset from the stored OOXML width by ~0.83 chars.)
$ws.Columns.Item(3).ColumnWidth = 53.17   # C: 56 -> 54
$ws.Columns.Item(4).ColumnWidth = 39.17   # D: 39 -> 40
$ws.Columns.Item(8).ColumnWidth = 35.17   # H: 29 -> 36
